$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear (contents + formatting) the PaymentType/FEINSSN "Y" marker cells in column C
# for every row except 19-24 (those rows get their execution timestamp refreshed instead).
$clearRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48)
foreach ($r in $clearRows) {
    $ws.Range("C$r").Clear()
}

# Refresh the execution timestamps recorded for the "Extension Payments" RAD run (rows 19-24)
$ws.Range("B19").Value = "Wed Mar 20 23:02:42 EDT 2024"
$ws.Range("B20").Value = "Wed Mar 20 23:02:52 EDT 2024"
$ws.Range("B21").Value = "Wed Mar 20 23:03:01 EDT 2024"
$ws.Range("B22").Value = "Wed Mar 20 23:03:11 EDT 2024"
$ws.Range("B23").Value = "Wed Mar 20 23:03:21 EDT 2024"
$ws.Range("B24").Value = "Wed Mar 20 23:03:30 EDT 2024"

# Update the active selection to mirror where the author left the cursor after the run.
$ws.Range("C41").Select()
